$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value to a literal piece of text without letting
# Excel auto-convert date-like strings (e.g. "01/01/2022") into date
# serial numbers (which would also spawn a brand-new cell style). We do
# this by entering the text as a formula returning a string literal and
# then "Paste Special -> Values" over itself, which collapses the
# formula down to a plain shared-string cell while keeping the
# original style untouched.
function Set-PlainText {
    param($range, [string]$text)
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163) | Out-Null
    $excel.CutCopyMode = $false
}

# --- Row 8: "Ativação:" date, 01/01/2016 -> 01/01/2022 ---
Set-PlainText $ws.Range("B8") "01/01/2022"
Set-PlainText $ws.Range("C8") "01/01/2022"

# --- Row 10: "Objetivos:" text rewritten ---
$objetivos = "A disciplina busca introduzir o aluno ao ambiente de engenharia, propondo problemas desafiadores gerando aptidão para solução de problemas. Apresentar a Engenharia de Materiais e seus campos de atuação, aspectos legais e éticos, bem como o mercado de trabalho para o engenheiro de materiais no Século XXI. Propiciar aos alunos uma visão geral do curso, com apresentação do currículo do curso de Engenharia de Materiais da EEL. Apresentar aos alunos uma visão da evolução histórica dos materiais com o homem. Descrever exemplos marcantes da introdução de novos materiais e as mudanças sociais provocadas. Apresentar o caráter interdisciplinar da Ciência e Engenharia de Materiais e suas ligações com outros ramos da Ciência. Apresentar estudos de caso demonstrando este caráter interdisciplinar."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# --- Insert new row 14 for the second responsible teacher ---
# Everything from the old row 14 ("Programa resumido:") downward shifts
# down by one row. Row.Insert() on this engine already mirrors the
# formatting (cell styles) of the row being pushed down, so no extra
# style bookkeeping is needed.
$ws.Rows.Item(14).Insert()
$ws.Range("B14").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C14").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"

# --- Row 15 (was 14): "Programa resumido:" text rewritten + height 60 ---
$programaResumido = "1- A importância dos materiais na evolução do homem na pré-história. Alquimia, Revolução Científica e a Revolução Industrial. 2-O Engenheiro como um profissional, funções da engenharia, a ética e comunicação na engenharia 3-A grandes áreas da Engenharia de Materiais. A interdisciplinaridade da Ciência e Engenharia de Materiais. 4- Perspectivas para a Engenharia de Materiais no século XXI. 5- O currículo do curso de engenharia de materiais da EEL-USP. 6- Noções básicas de Projetos em Engenharia.Em todos o conteúdo do curso serão abordados aspectos sociais, ambientais, éticos, legais e econômicos para ampliar as competências dos alunos"
$ws.Range("B15").Value = $programaResumido
$ws.Range("C15").Value = $programaResumido
$ws.Rows.Item(15).RowHeight = 60

# --- Row 16 (was 15): "Short syllabus:" height changes from 120 (shared
#     with row 17) down to its own 60 ---
$ws.Rows.Item(16).RowHeight = 60

# --- Row 17 (was 16): "Programa:" text rewritten (height stays 120) ---
$programa = "1- As características importantes de um engenheiro: aptidões interpessoais, aptidões de comunicação, liderança e competência. O engenheiro, profissional que busca solucionar problemas. 2-A Engenharia de Materiais: áreas de atuação e mercado de trabalho. Aplicação. A importância dos materiais na evolução do homem, as grandes áreas e interdisciplinaridade da Ciência e Engenharia de Materiais. Visita ao Departamento de Engenharia de Materiais. Conhecimento dos Grupos de Pesquisa do Departamento. Perspectivas para a Engenharia de Materiais no século XXI. 3- O campo de trabalho do engenheiro de materiais e suas áreas de atuação. Visita externa para integralização dos conhecimentos. 4- O currículo do curso de engenharia de materiais na EEL/USP. 5- Apresentação do método de trabalho com projetos, definindo os atributos de um projeto de engenharia, mapas conceituais e ferramentas que ilustram ideias e relações entre elas. Formular estratégias para resolução de problemas de engenharia. Estudo de casos"
$ws.Range("B17").Value = $programa
$ws.Range("C17").Value = $programa

# --- Row 18 (was 17): "Syllabus:" keeps height 120 (unchanged) ---

# --- Row 19 (was 18): "Avaliação:" never had a custom height, and the
#     row-insert above already carried that "no custom height" forward,
#     so nothing to do here. ---

# --- Row 20 (was 19): "Método:" text rewritten (height stays 60) ---
$metodo = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras"
$ws.Range("B20").Value = $metodo
$ws.Range("C20").Value = $metodo

# --- Row 21 (was 20): "Critério:" text rewritten (height stays 60) ---
$criterio = "Média Aritmética dos Projetos, Trabalhos, Relatórios e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."
$ws.Range("B21").Value = $criterio
$ws.Range("C21").Value = $criterio

# --- Row 22 (was 21): "Norma de recuperação:" text rewritten + height
#     shrinks from 120 (shared with row 22/bibliography) to its own 60 ---
$norma = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
$ws.Range("B22").Value = $norma
$ws.Range("C22").Value = $norma
$ws.Rows.Item(22).RowHeight = 60

# --- Row 23 (was 22): "Bibliografia:" text rewritten (height stays 120) ---
$biblio = "1) BROCKMAN, J.B. Introdução à Engenharia: Modelagem e Solução de Problemas, LTC Livros Científicos Editora, 2010.2) M.T. HOLTZAPPLE, W.D. REECE, Introdução à Engenharia: Modelagem e Solução de Problemas, LTC Livros Científicos Editora, 2006.2) CALLISTER Jr., W.D. Ciência e Engenharia de Materiais: Uma Introdução. LTC Livros Científicos Editora, 7a.ed., 2008. 4) - COHEN, M. (Ed.). Ciência e Engenharia de Materiais: sua Evolução, Prática e Perspectivas. Parte I: Materiais na história e na sociedade, 98p. Parte II: A Ciência e Engenharia de Materiais como uma multidisciplina, Tradução: José Roberto Gonçalves da Silva, São Carlos, UFSCar, 1985.5) Artigos científicos"
$ws.Range("B23").Value = $biblio
$ws.Range("C23").Value = $biblio
